# "Windows Azure" -> "Microsoft Azure" across the deck.
# Each target run is split into 3 runs ("...", "Microsoft ", "Azure...")
# so the new "Microsoft " piece lands in its own <a:r>, matching how
# PowerPoint itself fragments a run when only part of it is retyped.

$p = $ppt.ActivePresentation

# --- Slide 3: Agenda bullet "Deploying to Windows Azure" (54pt title run) ---
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$para3 = $tr3.Paragraphs(4, 1)
$para3.Text = "Deploying to "
$r3b = $para3.InsertAfter("Microsoft ")
$r3c = $r3b.InsertAfter("Azure")

# --- Slide 9: Title "Deploying to Windows Azure" ---
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(1).TextFrame.TextRange
$tr9.Text = "Deploying to "
$r9b = $tr9.InsertAfter("Microsoft ")
$r9c = $r9b.InsertAfter("Azure")

# --- Slide 10: Subtitle "Deploying to Windows Azure Web Sites" ---
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(2).TextFrame.TextRange
$tr10.Text = "Deploying to "
$r10b = $tr10.InsertAfter("Microsoft ")
$r10c = $r10b.InsertAfter("Azure Web Sites")

# --- Slide 11: Recap bullet "...to Windows Azure" (only last run changes) ---
$s11 = $p.Slides.Item(11)
$tr11 = $s11.Shapes.Item(2).TextFrame.TextRange
$full11 = $tr11.Text
$winIdx = $full11.IndexOf("Windows Azure")
$winStart = $winIdx + 1
$sub11 = $tr11.Characters($winStart, 8)
$sub11.Text = "Microsoft "
